$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178327441215515
$ws.Range("B1").Value = 1.872254014015198
$ws.Range("D1").Value = 0.4559695720672607
$ws.Range("E1").Value = 0.4881070256233215
